# Update Level/Editor task descriptions to reflect progress ("DONE" markers).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# D2: "Editor, Reactive" -> "Editor, Reactive. DONE"
$ws.Range("D2").Value = "Editor, Reactive. DONE"

# D9: "Editor: Camera, Reload, Terxture, Entities" -> "Editor: Camera DONE, Reload DONE, Terxture, Entities"
$ws.Range("D9").Value = "Editor: Camera DONE, Reload DONE, Terxture, Entities"
